$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4 data
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 4

# Update the active selection to match the target state (B4)
$ws.Range("B4").Select()
